# Adding Spreadsheet Writing Functionality
# The program now writes the important data findings (full text location,
# words sampled, collegiate-word count and ratio) into the "Sample Data"
# worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for a new leading column that records where the full text
# being sampled lives. This shifts the existing "Sample Location" ...
# "Collegiate Word Ratio" headers (and the sample row beneath them) one
# column to the right, from A:F to B:G.
$ws.Range("A1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "Full Text Location"

# The actual findings produced by the analysis run.
$ws.Range("E2").Value = 793
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 0.06305170239596469
